$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 is a new sales-order row appended below the existing data
# (A18:L18). Column A holds a date formatted/typed as plain text in this
# sheet (see A2/A14/A16), so force text formatting before assigning the
# value to stop Excel from auto-converting "02/14/2024" into a date serial.
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "02/14/2024"

$ws.Range("B18").Value = "SO240214001"
$ws.Range("C18").Value = "AB"
$ws.Range("D18").Value = "9999999999"
$ws.Range("E18").Value = "ARTIST"
$ws.Range("F18").Value = "TITLE"
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = "AEC"
$ws.Range("J18").Value = "LP"
$ws.Range("K18").Value = "Ashley"
$ws.Range("L18").Value = "No"
